$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the "No" index column (column A), shifting the hobby names (column B) into column A.
$ws.Columns.Item(1).Delete()

# Select the full column A, similar to the post-edit workbook state.
$ws.Columns.Item(1).Select()

$wb.Save()
